# Update cryptos list data (price & volume), and swap Maker/TrustWalletToken rows 41-42
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.298.25"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.858.38"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7009"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.24"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07868"
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  +5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08168"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.870.64"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.212"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7063"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.53"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "29.354.57"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.793"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007821"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.29"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "2.127.38"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.557"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.22"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.892"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.906"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.301"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.031"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05182"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.179"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7100"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9998"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.685"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.143.66"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9223"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.947"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4243"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.32"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.71"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5316"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.743"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.189"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.001"
$ws.Range("E51").Value = "  +0.67%  "
